$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the UML diagram text for the AddMovie method (D9) and the
# Movie array field declaration (D7) on Sheet1.
$ws.Range("D9").Value = " +AddMovie(Movie: newMovie&): bool "
$ws.Range("D7").Value = " -Movie: Movies[24]"

# Update the active selection to match the author's final cursor position.
$ws.Range("D7").Select()
